$wb = $excel.ActiveWorkbook

# Update header labels on the "FM" (Female Masters final) and "FF" (Female Final)
# worksheets: the Minute/Second/Rep columns lose their workout-specific suffix.
$wsFM = $wb.Worksheets.Item("FM")
$wsFF = $wb.Worksheets.Item("FF")

$wsFM.Range("D1").Value = "Minute"
$wsFM.Range("E1").Value = "Second"
$wsFM.Range("F1").Value = "Rep"

$wsFF.Range("D1").Value = "Minute"
$wsFF.Range("E1").Value = "Second"
$wsFF.Range("F1").Value = "Rep"

# Update the selection shown on each worksheet and flip which tab is active:
# FF was the active/selected tab, now FM is.
$wsFF.Range("D1:F1").Select()

$wsFM.Activate()
$wsFM.Range("D1:F1").Select()
